$d = $word.ActiveDocument

# Locate the paragraph that ends the existing bullet list:
# "Top 10 de los clientes que más paquetes han enviado"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Top 10 de los clientes que más paquetes han enviado*") {
        $target = $p
    }
}

$newTexts = @(
    "Promedio de clientes L-VIP por oficina (¿)",
    "Promedio anual de paquetes enviados de cada oficina",
    "Día de la semana con más envíos de paquetes ",
    "Top 5 de oficinas con menos ingresos",
    ""
)

$cur = $target
foreach ($t in $newTexts) {
    $cur.Range.InsertParagraphAfter()
    $cur = $cur.Next()
    if ($t -ne "") {
        $cur.Range.Text = $t
    }
}
